# Add "Topsec" as a new bold "Connections" entry, right after "APT1" and
# before the trailing blank paragraph. The trailing "_GoBack" bookmark
# (Word's auto-maintained "last edit" marker) should end up collapsed
# right after the new "Topsec" text, matching where Word would leave it
# after the last edit.

$d = $word.ActiveDocument

# Locate the "APT1" paragraph (the last entry currently in the
# "Connections:" list), so we don't depend on hard-coded paragraph
# indices.
$apt1Para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd("`r") -eq "APT1") {
        $apt1Para = $candidate
    }
}

# Insert a new empty paragraph right after "APT1".
$apt1Para.Range.InsertParagraphAfter()

# Re-fetch paragraphs; the new (currently empty) paragraph is the one
# immediately after "APT1".
$newIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd("`r") -eq "APT1") {
        $newIndex = $i + 1
    }
}
$newPara = $d.Paragraphs.Item($newIndex)

# Give it the same bold-only run formatting as its siblings, and fill in
# the text. A trailing sentinel character "X" is used temporarily: it
# gives the Bookmarks.Add call below a real, non-empty run of text to
# anchor on, and it is deleted again afterwards, leaving the bookmark
# collapsed right after "Topsec" without splitting the run.
$newPara.Range.Text = "TopsecX"
$newPara.Range.Bold = 1

$sentinelStart = $newPara.Range.Start + 6
$sentinelEnd = $newPara.Range.Start + 7
$sentinelRange = $d.Range($sentinelStart, $sentinelEnd)
$d.Bookmarks.Add("_GoBack", $sentinelRange)

$sentinelRange2 = $d.Range($sentinelStart, $sentinelEnd)
$sentinelRange2.Delete()
